$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.528.26"
$ws.Range("E2").Value = "  +0.91%  "

# Row 3
$ws.Range("D3").Value = "2.480.74"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.98"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.33"

# Row 7
$ws.Range("E7").Value = "  -0.96%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  +2.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.79"
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +1.34%  "

# Row 12
$ws.Range("E12").Value = "  +2.48%  "

# Row 13
$ws.Range("D13").Value = "2.861.97"
$ws.Range("E13").Value = "  +0.77%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.27"
$ws.Range("E14").Value = "  +9.97%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.87"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("D16").Value = "2.493.17"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.767"
$ws.Range("E17").Value = "  -1.77%  "

# Row 18
$ws.Range("D18").Value = "41.526.04"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("E19").Value = "  +2.49%  "

# Row 20
$ws.Range("E20").Value = "  +2.35%  "

# Row 21
$ws.Range("E21").Value = "  +4.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("E22").Value = "  +1.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.21"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  -1.37%  "

# Row 25
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.97"
$ws.Range("E27").Value = "  +4.32%  "

# Row 28
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("E29").Value = "  +1.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.19"
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.31"
$ws.Range("E31").Value = "  +3.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.47"
$ws.Range("E32").Value = "  +0.14%  "

# Row 33
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("E34").Value = "  +2.26%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.41"
$ws.Range("E35").Value = "  +3.16%  "

# Row 36
$ws.Range("E36").Value = "  -8.15%  "

# Row 37
$ws.Range("E37").Value = "  +4.64%  "

# Row 38
$ws.Range("E38").Value = "  -2.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  -2.17%  "

# Row 40
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("E41").Value = "  -1.64%  "

# Row 42
$ws.Range("E42").Value = "  -0.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.49"
$ws.Range("E43").Value = "  -2.79%  "

# Row 44
$ws.Range("D44").Value = "1.972.50"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  -1.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.97"
$ws.Range("E47").Value = "  +3.48%  "

# Row 48
$ws.Range("D48").Value = "2.719.94"
$ws.Range("E48").Value = "  +0.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.74"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.15"
$ws.Range("E50").Value = "  -1.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.66"
$ws.Range("E51").Value = "  -1.47%  "
